$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 41: now XGBoost / C Si / 96.8
$ws.Range("A41").Value = "XGBoost"
$ws.Range("B41").Value = "C Si"
$ws.Range("C41").Value = 96.8

# Update row 42: now Random Forest Classifier / C Si N / 96.9
$ws.Range("A42").Value = "Random Forest Classifier"
$ws.Range("B42").Value = "C Si N"
$ws.Range("C42").Value = 96.9

# Update row 43: now Random Forest Classifier / C Si / 97.3
$ws.Range("A43").Value = "Random Forest Classifier"
$ws.Range("B43").Value = "C Si"
$ws.Range("C43").Value = 97.3

# Move the active selection to D1
$ws.Range("D1").Select()
